$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 13 (shifts old rows 13-23 down to 14-24).
#    This creates the new "Docentes responsaveis:" value row (B13/C13).
$ws.Rows("13:13").Insert()

$ws.Range("B10").Value = 'Estudo das classes de compostos orgânicos mais relevantes do ponto de vista da sua presença no meio ambiente. Compreender a relação das estruturas moleculares às propriedades físicas e à sua reatividade química. Propor ações adequadas que possam eliminar ou mitigar efeitos danosos a elas relacionados a fim de possibilitar um crescimento e desenvolvimento sustentável para as futuras gerações.'
$ws.Range("C10").Value = 'Estudo das classes de compostos orgânicos mais relevantes do ponto de vista da sua presença no meio ambiente. Compreender a relação das estruturas moleculares às propriedades físicas e à sua reatividade química. Propor ações adequadas que possam eliminar ou mitigar efeitos danosos a elas relacionados a fim de possibilitar um crescimento e desenvolvimento sustentável para as futuras gerações.'
$ws.Range("B13").Value = '2346890 - Eliane Corrêa Pedrozo'
$ws.Range("C13").Value = '2346890 - Eliane Corrêa Pedrozo'
$ws.Range("B14").Value = 'Propriedades físicas, acidez e basicidade dos compostos orgânicos. Classes de compostos orgânicos, estereoquimica, principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Importância Industrial, impacto ambiental e degradação de compostos orgânicos.'
$ws.Range("C14").Value = 'Propriedades físicas, acidez e basicidade dos compostos orgânicos. Classes de compostos orgânicos, estereoquimica, principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Importância Industrial, impacto ambiental e degradação de compostos orgânicos.'
$ws.Range("B16").Value = 'Terminologia e conceitos básicos de química orgânica, estrutura e características das importantes classes dos compostos orgânicos, relacionando com produtos orgânicos antropogênicos. Acidez e basicidade em compostos orgânicos. Discussão dos aspectos básicos dos processos de partição. Como a estrutura química determina a solubilidade do composto orgânico em água. Estereoquimica. Principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Estereoquímica. Como os absorventes naturais são importantes para transporte, distribuição e destino dos compostos orgânicos no meio ambiente. Processos de transformação abiótica e biológicas. Processos químicos, fotoquímicos, biológicos e reações de degradação. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aquíferos.'
$ws.Range("C16").Value = 'Terminologia e conceitos básicos de química orgânica, estrutura e características das importantes classes dos compostos orgânicos, relacionando com produtos orgânicos antropogênicos. Acidez e basicidade em compostos orgânicos. Discussão dos aspectos básicos dos processos de partição. Como a estrutura química determina a solubilidade do composto orgânico em água. Estereoquimica. Principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Estereoquímica. Como os absorventes naturais são importantes para transporte, distribuição e destino dos compostos orgânicos no meio ambiente. Processos de transformação abiótica e biológicas. Processos químicos, fotoquímicos, biológicos e reações de degradação. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aquíferos.'
$ws.Range("B19").Value = 'Aulas expositivas sobre os tópicos descritos no programa. Seminários e/ou projeto dirigido. Provas escritas.'
$ws.Range("C19").Value = 'Aulas expositivas sobre os tópicos descritos no programa. Seminários e/ou projeto dirigido. Provas escritas.'
$ws.Range("B20").Value = 'Avaliação individual escrita e apresentação de seminário e/ou projeto. Avaliações individuais correspondem a 70% da NF (nota final) e em grupo 30% da NF. Se NF for => 5,0: Aluno aprovado, condicionado à frequência mínima de 70%.'
$ws.Range("C20").Value = 'Avaliação individual escrita e apresentação de seminário e/ou projeto. Avaliações individuais correspondem a 70% da NF (nota final) e em grupo 30% da NF. Se NF for => 5,0: Aluno aprovado, condicionado à frequência mínima de 70%.'
$ws.Range("B21").Value = 'Prova escrita envolvendo todo o conteúdo da disciplina. MF = (NF+R) / 2 & => 5,0 Aprovado'
$ws.Range("C21").Value = 'Prova escrita envolvendo todo o conteúdo da disciplina. MF = (NF+R) / 2 & => 5,0 Aprovado'
$ws.Range("B22").Value = '- Bruice, Paula Yurkanis – Química Orgânica – PEARSON Prentice Hall – São Paulo, 2006. Solomons, T.W.G.- Química Orgânica. Volumes 1 e 2, Rio de Janeiro, Livros Técnicos e Científicos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry – Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'
$ws.Range("C22").Value = '- Bruice, Paula Yurkanis – Química Orgânica – PEARSON Prentice Hall – São Paulo, 2006. Solomons, T.W.G.- Química Orgânica. Volumes 1 e 2, Rio de Janeiro, Livros Técnicos e Científicos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry – Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'

# 2. Row 13 keeps only B13/C13 (the new label "Docentes responsaveis:" value);
#    clear the stray A13 cell that Insert() left behind with inherited formatting
#    (matches the diff: row 13 has no <c r="A13"> at all).
$ws.Range("A13").Clear()

# 3. Re-apply the B/C column formatting (wrap text, vertical top alignment,
#    normal/red font) to the newly written cells so they keep the same look
#    as the rest of the sheet (style index 2 for column B, 3 for column C).
foreach ($r in @(10,13,14,16,19,20,21,22)) {
  $bcell = $ws.Range("B$r")
  $bcell.Font.Bold = $false
  $bcell.Font.Color = 0
  $bcell.WrapText = $true
  $bcell.VerticalAlignment = -4160

  $ccell = $ws.Range("C$r")
  $ccell.Font.Bold = $false
  $ccell.Font.Color = 255
  $ccell.WrapText = $true
  $ccell.VerticalAlignment = -4160
}

# 4. Split the merged column A/B width definition into its own <col> entry
#    for column A only (column B already has its own explicit width).
$ws.Columns("A").ColumnWidth = 29.75
